$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$sub = $wb.Worksheets.Item("Substrates and products")
$biomass = $wb.Worksheets.Item("Biomass")

# The "RR =" label that used to sit next to the CH4 row moves down to sit
# right above the (now relocated) "Reaction rate vector RR" block.
$sub.Range("N12").Clear()

# Make room for the biomass rows (xa, ..., xa_aer, xa_bac) plus three spacer
# rows above the "Reaction rate vector RR" block, which shifts down from
# row 13 to row 20.
$sub.Rows("13:19").Insert(-4121)   # xlShiftDown

# --- Row 13: xa ---------------------------------------------------------
$sub.Range("A13").Value = "xa"
$sub.Range("I13").Value = "Y"
$sub.Range("N13").Clear()
$sub.Range("N13").Style = "Normal"
$sub.Range("U13").Value = "xa"
$sub.Range("U13").Style = "Normal"
$sub.Range("Y13").Value = "xa"
$sub.Range("AC13").Value = "xa"
$sub.Range("AG13").Value = "xa"
$sub.Range("AG13").Style = "Normal"

# --- Row 14: ... ---------------------------------------------------------
$sub.Range("A14").Value = "..."
$sub.Range("I14").Value = "Y"
$sub.Range("K9").Copy()
$sub.Range("K14").PasteSpecial(-4122)   # xlPasteFormats
$sub.Range("K14").ClearContents()
$sub.Range("N14").Clear()
$sub.Range("N14").Style = "Normal"
$sub.Range("U14").Value = "..."
$sub.Range("U14").Style = "Normal"
$sub.Range("Y14").Value = "..."
$sub.Range("AC14").Value = "..."
$sub.Range("AG14").Value = "..."
$sub.Range("AG14").Style = "Normal"

# --- Row 15: xa_aer -------------------------------------------------------
$sub.Range("A15").Value = "xa_aer"
$sub.Range("B15:I15").Value = "Y"
$sub.Range("N15").Clear()
$sub.Range("U15").Value = "xa_aer"
$sub.Range("U15").Style = "Normal"
$sub.Range("Y15").Value = "xa_aer"
$sub.Range("AC15").Value = "xa_aer"
$sub.Range("AG15").Value = "xa_aer"
$sub.Range("AG15").Style = "Normal"

# --- Row 16: xa_bac -------------------------------------------------------
$sub.Range("A16").Value = "xa_bac"
$sub.Range("B16:I16").Value = "Y"
$sub.Range("N16").Clear()
$sub.Range("U16").Value = "xa_bac"
$sub.Range("U16").Style = "Normal"
$sub.Range("Y16").Value = "xa_bac"
$sub.Range("AC16").Value = "xa_bac"
$sub.Range("AG16").Value = "xa_bac"
$sub.Range("AG16").Style = "Normal"

# --- Rows 17-18: blank spacer rows ----------------------------------------
$sub.Range("I17").Clear()
$sub.Range("N17").Clear()
$sub.Range("AG17").Clear()
$sub.Range("U17").Style = "Normal"

$sub.Range("I18").Clear()
$sub.Range("N18").Clear()
$sub.Range("AG18").Clear()
$sub.Range("U18").Style = "Normal"

# --- Row 19: "RR =" label moves down here ---------------------------------
$sub.Range("I19").Clear()
$sub.Range("AG19").Clear()
$sub.Range("U19").Style = "Normal"
$sub.Range("N19").Value = "RR ="

# Remove the now-redundant Biomass sheet; its rows were folded into
# "Substrates and products" above.
$biomass.Delete()

$sub.Range("AH11").Select()
